$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 71432130
$ws.Range("I86").Value = 500000000
$ws.Range("J86").Value = 4150
$ws.Range("K86").Value = 500000000
$ws.Range("L86").Value = 4150
$ws.Range("M86").Value = -499998877
$ws.Range("N86").Value = -6396
$ws.Range("H89").Value = 71432130
$ws.Range("I89").Value = 500000000
$ws.Range("J89").Value = 4150
$ws.Range("K89").Value = 2500000000
$ws.Range("L89").Value = 20750
$ws.Range("M89").Value = -2499994384
$ws.Range("N89").Value = -31982
$ws.Range("H113").Value = 2833.3333
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -11008
$ws.Range("H129").Value = 1611.8036
$ws.Range("I129").Value = 382.81818
$ws.Range("J129").Value = 1912.2222
$ws.Range("K129").Value = 1148.45454
$ws.Range("L129").Value = 5736.6666
$ws.Range("M129").Value = 3851.54546
$ws.Range("N129").Value = -15736.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1562.6364
$ws.Range("I2").Value = 1139.4286
$ws.Range("J2").Value = 2303.25
$ws.Range("K2").Value = 1139.4286
$ws.Range("L2").Value = 2303.25
$ws.Range("M2").Value = -1026.4286
$ws.Range("N2").Value = -2529.25
$ws.Range("H32").Value = 9043.312
$ws.Range("I32").Value = 9320.725
$ws.Range("K32").Value = 9320.725
$ws.Range("M32").Value = -9033.725
$ws.Range("H45").Value = 1446.7894
$ws.Range("I45").Value = 1277.8572
$ws.Range("J45").Value = 1919.8
$ws.Range("K45").Value = 1277.8572
$ws.Range("L45").Value = 1919.8
$ws.Range("M45").Value = -900.8571999999999
$ws.Range("N45").Value = -2673.8
$ws.Range("H110").Value = 2225
$ws.Range("I110").Value = 2128.5715
$ws.Range("K110").Value = 2128.5715
$ws.Range("M110").Value = -83.57150000000001
$ws.Range("H116").Value = 1562.6364
$ws.Range("I116").Value = 1139.4286
$ws.Range("J116").Value = 2303.25
$ws.Range("K116").Value = 1139.4286
$ws.Range("L116").Value = 2303.25
$ws.Range("M116").Value = 1154.5714
$ws.Range("N116").Value = -6891.25
$ws.Range("H122").Value = 6884.65
$ws.Range("I122").Value = 7225.9473
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 21677.8419
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = -19227.8419
$ws.Range("N122").Value = -6100
$ws.Range("H132").Value = 10871469
$ws.Range("I132").Value = 13890175
$ws.Range("K132").Value = 41670525
$ws.Range("M132").Value = -41667995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1562.6364
$ws.Range("I3").Value = 1139.4286
$ws.Range("J3").Value = 2303.25
$ws.Range("K3").Value = 1139.4286
$ws.Range("L3").Value = 2303.25
$ws.Range("M3").Value = -1025.4286
$ws.Range("N3").Value = -2531.25
$ws.Range("H26").Value = 25607.6
$ws.Range("I26").Value = 9000
$ws.Range("K26").Value = 9000
$ws.Range("M26").Value = -8708
$ws.Range("H107").Value = 2166.6667
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -80
$ws.Range("N107").Value = -6340
$ws.Range("H134").Value = 3066.7874
$ws.Range("I134").Value = 1968.8387
$ws.Range("K134").Value = 5906.5161
$ws.Range("M134").Value = -3371.5161

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2080
$ws.Range("I16").Value = 2080
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2080
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1793
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 10758342
$ws.Range("I31").Value = 11171
$ws.Range("J31").Value = 16669286
$ws.Range("K31").Value = 11171
$ws.Range("L31").Value = 16669286
$ws.Range("M31").Value = -10876
$ws.Range("N31").Value = -16669876
$ws.Range("H34").Value = 10758342
$ws.Range("I34").Value = 11171
$ws.Range("J34").Value = 16669286
$ws.Range("K34").Value = 11171
$ws.Range("L34").Value = 16669286
$ws.Range("M34").Value = -10969
$ws.Range("N34").Value = -16669690
$ws.Range("H58").Value = 3175.3872
$ws.Range("I58").Value = 717.0952
$ws.Range("K58").Value = 717.0952
$ws.Range("M58").Value = -514.0952
$ws.Range("H99").Value = 1982.4
$ws.Range("I99").Value = 1878
$ws.Range("K99").Value = 1878
$ws.Range("M99").Value = -380
$ws.Range("H105").Value = 1655.7142
$ws.Range("I105").Value = 1520
$ws.Range("J105").Value = 1995
$ws.Range("K105").Value = 1520
$ws.Range("L105").Value = 1995
$ws.Range("M105").Value = 227
$ws.Range("N105").Value = -5489
$ws.Range("H113").Value = 2080
$ws.Range("I113").Value = 2080
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2080
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 90
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3363.6
$ws.Range("I122").Value = 3976
$ws.Range("J122").Value = 914
$ws.Range("K122").Value = 11928
$ws.Range("L122").Value = 2742
$ws.Range("M122").Value = -9478
$ws.Range("N122").Value = -7642
$ws.Range("H126").Value = 1982.4
$ws.Range("I126").Value = 1878
$ws.Range("K126").Value = 5634
$ws.Range("M126").Value = -3164
$ws.Range("H132").Value = 20003062
$ws.Range("I132").Value = 25002888
$ws.Range("J132").Value = 3762.4
$ws.Range("K132").Value = 75008664
$ws.Range("L132").Value = 11287.2
$ws.Range("M132").Value = -75006134
$ws.Range("N132").Value = -16347.2
$ws.Range("H134").Value = 1085541.4
$ws.Range("I134").Value = 3650.4375
$ws.Range("K134").Value = 10951.3125
$ws.Range("M134").Value = -8416.3125
$ws.Range("H136").Value = 3175.3872
$ws.Range("I136").Value = 717.0952
$ws.Range("K136").Value = 2151.2856
$ws.Range("M136").Value = 398.7143999999998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 143988.86
$ws.Range("I113").Value = 200921.8
$ws.Range("J113").Value = 1656.5
$ws.Range("K113").Value = 200921.8
$ws.Range("L113").Value = 1656.5
$ws.Range("M113").Value = -198751.8
$ws.Range("N113").Value = -5996.5
$ws.Range("H118").Value = 14190.488
$ws.Range("J118").Value = 14190.488
$ws.Range("L118").Value = 14190.488
$ws.Range("N118").Value = -17504.488
$ws.Range("H132").Value = 3639.76
$ws.Range("I132").Value = 3631
$ws.Range("J132").Value = 3655.3333
$ws.Range("K132").Value = 10893
$ws.Range("L132").Value = 10965.9999
$ws.Range("M132").Value = -8363
$ws.Range("N132").Value = -16025.9999
$ws.Range("H138").Value = 50421.5
$ws.Range("J138").Value = 50421.5
$ws.Range("L138").Value = 50421.5
$ws.Range("N138").Value = -60701.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1158.9286
$ws.Range("I16").Value = 1178.8462
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 1178.8462
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -1008.8462
$ws.Range("N16").Value = -1240
$ws.Range("H40").Value = 9393.333000000001
$ws.Range("I40").Value = 14500
$ws.Range("J40").Value = 7536.364
$ws.Range("K40").Value = 14500
$ws.Range("L40").Value = 7536.364
$ws.Range("M40").Value = -14364
$ws.Range("N40").Value = -7808.364
$ws.Range("H55").Value = 780
$ws.Range("I55").Value = 1000
$ws.Range("K55").Value = 1000
$ws.Range("M55").Value = -827
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -798
$ws.Range("N61").ClearContents()
$ws.Range("H94").Value = 50048.332
$ws.Range("J94").Value = 50048.332
$ws.Range("L94").Value = 50048.332
$ws.Range("N94").Value = -51400.332
$ws.Range("H110").Value = 40572
$ws.Range("J110").Value = 40572
$ws.Range("L110").Value = 40572
$ws.Range("N110").Value = -48752
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1170
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 6484.96
$ws.Range("I122").Value = 6741.6
$ws.Range("J122").Value = 6100
$ws.Range("K122").Value = 20224.8
$ws.Range("L122").Value = 18300
$ws.Range("M122").Value = -17774.8
$ws.Range("N122").Value = -23200
$ws.Range("H136").Value = 29413820
$ws.Range("I136").Value = 38463748
$ws.Range("J136").Value = 1548.5
$ws.Range("K136").Value = 115391244
$ws.Range("L136").Value = 4645.5
$ws.Range("M136").Value = -115388694
$ws.Range("N136").Value = -9745.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2583.3333
$ws.Range("I113").Value = 300
$ws.Range("J113").Value = 3725
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 11175
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -15515
$ws.Range("H122").Value = 2181.6428
$ws.Range("I122").Value = 2194.85
$ws.Range("J122").Value = 2148.625
$ws.Range("K122").Value = 6584.549999999999
$ws.Range("L122").Value = 6445.875
$ws.Range("M122").Value = -4134.549999999999
$ws.Range("N122").Value = -11345.875
$ws.Range("H126").Value = 2826.348
$ws.Range("I126").Value = 2077.75
$ws.Range("J126").Value = 4537.4287
$ws.Range("K126").Value = 6233.25
$ws.Range("L126").Value = 13612.2861
$ws.Range("M126").Value = -3763.25
$ws.Range("N126").Value = -18552.2861
$ws.Range("H132").Value = 1137.8431
$ws.Range("I132").Value = 1023.39026
$ws.Range("J132").Value = 1607.1
$ws.Range("K132").Value = 3070.17078
$ws.Range("L132").Value = 4821.299999999999
$ws.Range("M132").Value = -540.1707799999999
$ws.Range("N132").Value = -9881.299999999999
